$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 63, shifting existing rows 63:83 down to 64:84
# (mirrors a native Excel "Insert Row" / xlShiftDown operation)
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with the new weekly price record
$ws.Range("A63").Value = 3
$ws.Range("B63").Value = "Femacal de La Calera"
$ws.Range("C63").Value = "Coquimbo"
$ws.Range("D63").Value = 44809
$ws.Range("E63").Value = 5
$ws.Range("F63").Value = 100112035
$ws.Range("G63").Value = "Bruselas (repollito)"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 105
$ws.Range("K63").Value = 15000
$ws.Range("L63").Value = 16000
$ws.Range("M63").Value = 15476
$ws.Range("N63").Value = "`$/malla 15 kilos"
$ws.Range("O63").Value = "Provincia de Quillota"
$ws.Range("P63").Value = 1032
$ws.Range("Q63").Value = 15
$ws.Range("R63").Value = "Hortaliza"
